$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'300.54"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'0.52%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'32.16"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'1.91%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'4.972"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-3.55%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07878"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-2.70%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'2.085"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-15.78%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'7.793"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-0.04%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'3.836"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-2.03%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.9258"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-0.39%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1739"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-1.30%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07987"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'7.59%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.08688"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-1.83%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.03099"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'2.95%"
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'0.20%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001512"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-1.12%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.005855"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-2.66%"
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'2,098.99%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.463"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-1.85%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'2.256"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-1.48%"
$ws.Range("E19").Style = "Normal"
$ws.Range("E21").Value = "'-2.29%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'4.312"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'3.76%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.1795"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'6.75%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.04607"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-0.32%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.001238"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-0.27%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.004433"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-2.10%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0001252"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'4.35%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.01720"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-1.69%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.04750"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'3.33%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007493"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'8.40%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1356"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-1.18%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002364"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'7.96%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.01130"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'9.78%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005996"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-3.38%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00000000751"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'0.14%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.003395"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-59.60%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.8205"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'9.62%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002103"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'0.14%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0002003"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.14%"
$ws.Range("E50").Style = "Normal"
